$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert two new columns before column G (shifts G:Q -> I:S, R:S -> T:U)
$ws.Range("G1:H1").EntireColumn.Insert()

# 2. Set widths for the two newly inserted columns (closest achievable given
#    this runtime's column-width quantization grid)
$ws.Columns("G").ColumnWidth = 52
$ws.Columns("H").ColumnWidth = 42.285714285714285

# 3. The old "Ridge Regression + normalization" columns shifted from R:S to
#    T:U; that data set is being retired, so clear its content (column width
#    metadata is left untouched, matching the target state).
$ws.Range("T:U").ClearContents()

# 4. Fix the header text in D1 ("Accuracy((LinearRegression))" -> "Accuracy(LinearRegression)")
$ws.Range("D1").Value = "Accuracy(LinearRegression)"

# 5. Populate the headers for the two new columns
$ws.Range("G1").Value = "Predicted values(Polynomial Prediction+ normalization)"
$ws.Range("H1").Value = "Accuracy(Polynomial Prediction+ normalization)"

# 6. Fill in the new data for rows 2-37 (G = predicted, H = accuracy)
$data = @(
    @(2, 87.510650975185598, 97.777263659425302),
    @(3, 84.847247169212395, 93.753864275372806),
    @(4, 80.000087429361102, 86.021599386409804),
    @(5, 82.462586236712596, 86.802722354434394),
    @(6, 85.611473552453106, 87.8066395409775),
    @(7, 95.342039535381005, 93.933043877222602),
    @(8, 99.964272048467905, 96.583837727988296),
    @(9, 100.726332090222, 94.136758962824501),
    @(10, 105.406160451524, 97.598296714374499),
    @(11, 111.39783552433499, 97.800150895104693),
    @(12, 106.77871177753801, 95.338135515659403),
    @(13, 110.844786606558, 97.232268953121505),
    @(14, 106.734418242817, 92.812537602450007),
    @(15, 101.316867636268, 88.101624031538194),
    @(16, 104.587709957659, 90.161818929016405),
    @(17, 102.812179413859, 87.873657618682898),
    @(18, 110.168588693232, 93.760501015516795),
    @(19, 105.722072597177, 88.842077812754198),
    @(20, 112.50215237002899, 93.751793641691094),
    @(21, 119.130899595331, 99.275749662775993),
    @(22, 119.589344759382, 99.657787299485506),
    @(23, 114.475128765305, 95.395940637754606),
    @(24, 110.394819500148, 91.9956829167905),
    @(25, 126.549281821127, 94.54226514906),
    @(26, 139.01200242618901, 84.156664644841996),
    @(27, 120.33632895205299, 99.719725873288795),
    @(28, 140.35254722858599, 84.956928501158501),
    @(29, 138.87508908097499, 89.781675332559402),
    @(30, 154.580143737962, 77.317346239712407),
    @(31, 148.788777069785, 81.913668992233994),
    @(32, 151.08733991988501, 80.089412761996002),
    @(33, 150.040824454155, 82.780605895190902),
    @(34, 148.478567873023, 84.001118849200694),
    @(35, 152.94317905253899, 80.513141365203197),
    @(36, 143.68963975029899, 87.742468945078301),
    @(37, 160.61577382394501, 74.518926700042499),
)

foreach ($row in $data) {
    $r = [int]$row[0]
    $ws.Cells.Item($r, 7).Value = $row[1]
    $ws.Cells.Item($r, 8).Value = $row[2]
}

# 7. Row 38 summary row: G38 gets the shared "Avg" label (plain, unstyled,
#    matching how Excel carries over the left neighbour's format on insert)
#    and H38 gets the new accuracy average figure.
$ws.Range("G38").Value = "Avg"
$ws.Range("H38").Value = 90.234658396692694

# 8. Final selection, matching the saved session state
$ws.Range("G39").Select()
